$wb = $excel.ActiveWorkbook

# 1. Rename sheet "good_too" -> "good_2"
$ws2 = $wb.Worksheets.Item("good_too")
$ws2.Name = "good_2"

# 2. Update the comment text on sheet "good", cell B1 - tighten the
#    worksheet-name validation rule wording.
$ws1 = $wb.Worksheets.Item("good")
$cmt = $ws1.Range("B1").Comment

$newText = "This worksheet is valid and will yield a METS <dmdSec> element if this file is passed to TOMES Packager.`n`n" +
    "It is valid because it contains the valid headers ""dc_element"" and ""dc_value"".`n`n" +
    "It is the user's responsibility to use correct Dublin Core element names in the ""dc_element"" column and legal XML values in the ""dc_value"" column.`n`n" +
    "The name of the worksheet (e.g. ""good"") will be used for the <dmdSec> element's ""ID"" attribute value. Use only letters, underscores, and numbers provided the name does not start with a number.`n`n" +
    "Multiple worksheets within the same Excel file may be used. They will each yield a new METS <dmdSec> element."

[void]$cmt.Text($newText)
